# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order a few country rows (shared-string ordering fix) ---
# Armenia / Austria / Moldavia block
$ws.Range("A59").Value = "Armenia"
$ws.Range("A60").Value = "Austria"
$ws.Range("A61").Value = "Moldavia"

# Hungria / Libia block
$ws.Range("A70").Value = "Hungria"
$ws.Range("A71").Value = "Libia"

# Georgia / Camerun block
$ws.Range("A92").Value = "Georgia"
$ws.Range("A93").Value = "Camerun"

# --- Update "last refreshed" timestamp caption ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 09:08"

# --- Refresh country statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 8584850
$ws.Range("C4").Value = 31
$ws.Range("D4").Value = 5602122
$ws.Range("E4").Value = 2755319

$ws.Range("B27").Value = 322879
$ws.Range("C27").Value = 7053
$ws.Range("D27").Value = 134898
$ws.Range("E27").Value = 181938
$ws.Range("G27").Value = 116
$ws.Range("H27").Value = 6043

$ws.Range("B28").Value = 307765
$ws.Range("C28").Value = 430
$ws.Range("D28").Value = 286119
$ws.Range("E28").Value = 19354
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 2292

$ws.Range("B59").Value = 70836
$ws.Range("C59").Value = 2306
$ws.Range("D59").Value = 49787
$ws.Range("E59").Value = 19918
$ws.Range("G59").Value = 10
$ws.Range("H59").Value = 1131

$ws.Range("B60").Value = 69409
$ws.Range("D60").Value = 52617
$ws.Range("E60").Value = 15867
$ws.Range("H60").Value = 925

$ws.Range("B61").Value = 68791
$ws.Range("D61").Value = 49702
$ws.Range("E61").Value = 17459
$ws.Range("H61").Value = 1630

$ws.Range("B63").Value = 64335
$ws.Range("C63").Value = 325
$ws.Range("E63").Value = 2374
$ws.Range("G63").Value = 4
$ws.Range("H63").Value = 538

$ws.Range("B70").Value = 52212
$ws.Range("C70").Value = 2032
$ws.Range("D70").Value = 15254
$ws.Range("E70").Value = 35653
$ws.Range("G70").Value = 46
$ws.Range("H70").Value = 1305

$ws.Range("B71").Value = 51625
$ws.Range("D71").Value = 28440
$ws.Range("E71").Value = 22420
$ws.Range("H71").Value = 765

$ws.Range("E85").Value = 3514
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 936

$ws.Range("B87").Value = 27466
$ws.Range("C87").Value = 22
$ws.Range("D87").Value = 25159
$ws.Range("E87").Value = 1402

$ws.Range("B92").Value = 22803
$ws.Range("C92").Value = 1595
$ws.Range("D92").Value = 9401
$ws.Range("E92").Value = 13224
$ws.Range("G92").Value = 6
$ws.Range("H92").Value = 178

$ws.Range("B93").Value = 21570
$ws.Range("D93").Value = 20117
$ws.Range("E93").Value = 1028
$ws.Range("H93").Value = 425

$ws.Range("B145").Value = 3897
$ws.Range("C145").Value = 7
$ws.Range("E145").Value = 629
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 102

$ws.Range("B176").Value = 583
$ws.Range("C176").Value = 2
$ws.Range("D176").Value = 545
$ws.Range("E176").Value = 31

$ws.Range("B203").Value = 39
$ws.Range("C203").Value = 1
$ws.Range("E203").Value = 12
